# "new file mobile checklist" - add the new "Conpatibility testing" section
# (section 2) to the bottom of the "General checklist" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Section header row: ID = 2, title = new shared string.
$ws.Range("A32").Value = 2
$ws.Range("B32").Value = "Conpatibility testing (тестування сумісності)"

# Two empty sub-items under the new section: 2.1 and 2.2.
$ws.Range("A33").Value = 2.1
$ws.Range("B33").Value = ""

$ws.Range("A34").Value = 2.2
$ws.Range("B34").Value = ""

# Move the selection to match the author's last selection before saving.
$ws.Range("C2:C31").Select()
